# Refresh the bundled test-data workbook:
#  1. Roll the timestamp embedded in the sample e-mail addresses on the
#     "UsuariosRegistro" sheet forward from 20251110_120951 to
#     20251110_130229.
#  2. Swap the sample catalogue rows on the "ProductosBusqueda" sheet for a
#     new set of sample products (laptops/phones/cameras/tablets instead of
#     desktops/monitors), and let column C re-fit to the new content.

$wb = $excel.ActiveWorkbook

# --- 1. Refresh e-mail timestamps everywhere they appear ------------------
# The same addresses are reused as valid-login fixtures on "LoginData", so
# both sheets need to move from the old run's timestamp to the new one.
$oldStamp = "20251110_120951"
$newStamp = "20251110_130229"

$wsUsuarios = $wb.Worksheets.Item("UsuariosRegistro")
for ($row = 2; $row -le 6; $row++) {
    $cell = $wsUsuarios.Cells.Item($row, 3)
    $oldValue = [string]$cell.Value2
    $cell.Value = $oldValue.Replace($oldStamp, $newStamp)
}

$wsLogin = $wb.Worksheets.Item("LoginData")
for ($row = 2; $row -le 3; $row++) {
    $cell = $wsLogin.Cells.Item($row, 1)
    $oldValue = [string]$cell.Value2
    $cell.Value = $oldValue.Replace($oldStamp, $newStamp)
}

# --- 2. Replace the sample product rows on ProductosBusqueda --------------
$wsProductos = $wb.Worksheets.Item("ProductosBusqueda")

# Columns: A = Categoria, B = SubCategoria, C = Producto, D = Cantidad
$wsProductos.Range("A2").Value = "Laptops & Notebooks"
$wsProductos.Range("B2").Value = "'"
$wsProductos.Range("C2").Value = "MacBook"

$wsProductos.Range("A3").Value = "'"
$wsProductos.Range("C3").Value = "iPhone"
$wsProductos.Range("D3").NumberFormat = "@"
$wsProductos.Range("D3").Value = "1"

$wsProductos.Range("A4").Value = "Cameras"
$wsProductos.Range("B4").Value = "'"
$wsProductos.Range("C4").Value = "Canon EOS 5D"

$wsProductos.Range("A5").Value = "Laptops & Notebooks"
$wsProductos.Range("B5").Value = "Macs"
$wsProductos.Range("C5").Value = "MacBook Air"
$wsProductos.Range("D5").NumberFormat = "@"
$wsProductos.Range("D5").Value = "2"

$wsProductos.Range("A6").Value = "Tablets"
$wsProductos.Range("C6").Value = "Samsung Galaxy Tab 10.1"

# The "Producto" column now holds longer names (e.g. "Samsung Galaxy Tab
# 10.1") -- re-fit its width to the new content, same as Excel would do
# when a user double-clicks the column border.
$wsProductos.Columns.Item(3).AutoFit() | Out-Null
